$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rodinia benchmarks (rows 60-76, 78-79) were run on CPU: populate "CPU time" (col K)
# and refresh "GPU time" (col L) with the new measured values. Both columns store
# their numeric readings as literal text (existing sheet convention), so we briefly
# force a Text number format before writing the value, then restore the original
# numeric display format so the cell style index is unchanged.
$numFmt = "#,##0.00"

$ws.Range("K60").NumberFormat = "@"
$ws.Range("K60").Value = "0.0823160000"
$ws.Range("K60").NumberFormat = $numFmt

$ws.Range("L60").NumberFormat = "@"
$ws.Range("L60").Value = "0.0975790000"
$ws.Range("L60").NumberFormat = $numFmt

$ws.Range("K61").NumberFormat = "@"
$ws.Range("K61").Value = "0.6391930000"
$ws.Range("K61").NumberFormat = $numFmt

$ws.Range("L61").NumberFormat = "@"
$ws.Range("L61").Value = "0.6694910000"
$ws.Range("L61").NumberFormat = $numFmt

$ws.Range("K62").NumberFormat = "@"
$ws.Range("K62").Value = "0.1200500000"
$ws.Range("K62").NumberFormat = $numFmt

$ws.Range("L62").NumberFormat = "@"
$ws.Range("L62").Value = "0.7782530000"
$ws.Range("L62").NumberFormat = $numFmt

$ws.Range("K63").NumberFormat = "@"
$ws.Range("K63").Value = "8.6193270000"
$ws.Range("K63").NumberFormat = $numFmt

$ws.Range("L63").NumberFormat = "@"
$ws.Range("L63").Value = "8.4539560000"
$ws.Range("L63").NumberFormat = $numFmt

$ws.Range("K64").NumberFormat = "@"
$ws.Range("K64").Value = "0.1033040000"
$ws.Range("K64").NumberFormat = $numFmt

$ws.Range("L64").NumberFormat = "@"
$ws.Range("L64").Value = "0.1483920000"
$ws.Range("L64").NumberFormat = $numFmt

$ws.Range("K65").NumberFormat = "@"
$ws.Range("K65").Value = "4.0266510000"
$ws.Range("K65").NumberFormat = $numFmt

$ws.Range("L65").NumberFormat = "@"
$ws.Range("L65").Value = "5.4771270000"
$ws.Range("L65").NumberFormat = $numFmt

$ws.Range("K66").NumberFormat = "@"
$ws.Range("K66").Value = "12.3171850000"
$ws.Range("K66").NumberFormat = $numFmt

$ws.Range("L66").NumberFormat = "@"
$ws.Range("L66").Value = "12.2726530000"
$ws.Range("L66").NumberFormat = $numFmt

$ws.Range("K67").NumberFormat = "@"
$ws.Range("K67").Value = "32.7613950000"
$ws.Range("K67").NumberFormat = $numFmt

$ws.Range("L67").NumberFormat = "@"
$ws.Range("L67").Value = "33.2539140000"
$ws.Range("L67").NumberFormat = $numFmt

$ws.Range("K68").NumberFormat = "@"
$ws.Range("K68").Value = "0.9711030000"
$ws.Range("K68").NumberFormat = $numFmt

$ws.Range("L68").NumberFormat = "@"
$ws.Range("L68").Value = "1.6896860000"
$ws.Range("L68").NumberFormat = $numFmt

$ws.Range("K69").NumberFormat = "@"
$ws.Range("K69").Value = "0.2438150000"
$ws.Range("K69").NumberFormat = $numFmt

$ws.Range("L69").NumberFormat = "@"
$ws.Range("L69").Value = "0.3204090000"
$ws.Range("L69").NumberFormat = $numFmt

$ws.Range("K70").NumberFormat = "@"
$ws.Range("K70").Value = "31.5730360000"
$ws.Range("K70").NumberFormat = $numFmt

$ws.Range("L70").NumberFormat = "@"
$ws.Range("L70").Value = "32.2186160000"
$ws.Range("L70").NumberFormat = $numFmt

$ws.Range("K71").NumberFormat = "@"
$ws.Range("K71").Value = "18.7550030000"
$ws.Range("K71").NumberFormat = $numFmt

$ws.Range("L71").NumberFormat = "@"
$ws.Range("L71").Value = "19.6663000000"
$ws.Range("L71").NumberFormat = $numFmt

$ws.Range("K72").NumberFormat = "@"
$ws.Range("K72").Value = "1.7172340000"
$ws.Range("K72").NumberFormat = $numFmt

$ws.Range("L72").NumberFormat = "@"
$ws.Range("L72").Value = "1.9030160000"
$ws.Range("L72").NumberFormat = $numFmt

$ws.Range("K73").NumberFormat = "@"
$ws.Range("K73").Value = "0.0606480000"
$ws.Range("K73").NumberFormat = $numFmt

$ws.Range("L73").NumberFormat = "@"
$ws.Range("L73").Value = "0.0606480000"
$ws.Range("L73").NumberFormat = $numFmt

$ws.Range("K74").NumberFormat = "@"
$ws.Range("K74").Value = "0.4908550000"
$ws.Range("K74").NumberFormat = $numFmt

$ws.Range("L74").NumberFormat = "@"
$ws.Range("L74").Value = "1.5245870000"
$ws.Range("L74").NumberFormat = $numFmt

$ws.Range("K75").NumberFormat = "@"
$ws.Range("K75").Value = "0.4632210000"
$ws.Range("K75").NumberFormat = $numFmt

$ws.Range("L75").NumberFormat = "@"
$ws.Range("L75").Value = "1.4709390000"
$ws.Range("L75").NumberFormat = $numFmt

$ws.Range("K76").NumberFormat = "@"
$ws.Range("K76").Value = "0.3353120000"
$ws.Range("K76").NumberFormat = $numFmt

$ws.Range("L76").NumberFormat = "@"
$ws.Range("L76").Value = "1.3506390000"
$ws.Range("L76").NumberFormat = $numFmt

$ws.Range("K78").NumberFormat = "@"
$ws.Range("K78").Value = "11.2296740000"
$ws.Range("K78").NumberFormat = $numFmt

$ws.Range("L78").NumberFormat = "@"
$ws.Range("L78").Value = "65.7484890000"
$ws.Range("L78").NumberFormat = $numFmt

$ws.Range("K79").NumberFormat = "@"
$ws.Range("K79").Value = "1.1714820000"
$ws.Range("K79").NumberFormat = $numFmt

$ws.Range("L79").NumberFormat = "@"
$ws.Range("L79").Value = "1.1548690000"
$ws.Range("L79").NumberFormat = $numFmt

